$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INS")

# --- Table 1 (rows 3-13): header renamed, sequential "Value" numbering,
#     and the separate "GAMS statement N" rows get merged into the
#     commented-out text of the preceding row (prefixed with "*").
$ws.Range("C3").Value  = "~TFM_INS"
$ws.Range("F4").Value  = "Value"

$ws.Range("F6").Value  = 2
$ws.Range("F7").Value  = 3

$ws.Range("E8").Value  = "*GAMS statement 1"
$ws.Range("F8").Value  = 4

$ws.Range("E9").Value  = "*GAMS statement 2"
$ws.Range("F9").Value  = 5

$ws.Range("E10").Value = "*GAMS statement 3"
$ws.Range("F10").Value = 6

$ws.Range("F11").Value = 7

$ws.Range("E12").Value = "*GAMS statement A"
$ws.Range("F12").Value = 8

$ws.Range("E13").Value = "*GAMS statement B"
$ws.Range("F13").Value = 9

# --- Table 2 (rows 17-33): same header rename + sequential numbering.
$ws.Range("C17").Value = "~TFM_INS"
$ws.Range("F18").Value = "Value"

$ws.Range("F20").Value = 2
$ws.Range("F21").Value = 3
$ws.Range("F22").Value = 4
$ws.Range("F23").Value = 5
$ws.Range("F24").Value = 6
$ws.Range("F25").Value = 7
$ws.Range("F26").Value = 8
$ws.Range("F27").Value = 9
$ws.Range("F28").Value = 10
$ws.Range("F29").Value = 11
$ws.Range("F30").Value = 12
$ws.Range("F31").Value = 13
$ws.Range("F32").Value = 14
$ws.Range("F33").Value = 15

# --- Table 3 (rows 37-53): same header rename + sequential numbering.
$ws.Range("C37").Value = "~TFM_INS"
$ws.Range("F38").Value = "Value"

$ws.Range("F40").Value = 2
$ws.Range("F41").Value = 3
$ws.Range("F42").Value = 4
$ws.Range("F43").Value = 5
$ws.Range("F44").Value = 6
$ws.Range("F45").Value = 7
$ws.Range("F46").Value = 8
$ws.Range("F47").Value = 9
$ws.Range("F48").Value = 10
$ws.Range("F49").Value = 11
$ws.Range("F50").Value = 12
$ws.Range("F51").Value = 13
$ws.Range("F52").Value = 14
$ws.Range("F53").Value = 15

# --- View state: scroll back to top-left default and move the selection.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E14").Select()
